# Update the "Förändrad" (Changed) date column (column C) from 2023-10-08
# (serial 45207) to 2023-10-09 (serial 45208) for data rows 2 through 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}
